# Update the "last added logger download" row on each site sheet with the
# newest readings (salt routes pulled 2020-09-22), and move the selection /
# active-tab bookmark from PBSF to WIC as it was left when the workbook was
# saved.

$wb = $excel.ActiveWorkbook

# --- PBSF: row 24 gets the new reading, selection moves to D30 ---------
$wsPBSF = $wb.Worksheets.Item("PBSF")
$wsPBSF.Range("A24").Value = 44096.468055555553
$wsPBSF.Range("B24").Value = 1566
$wsPBSF.Range("C24").Value = 20.399999999999999
$wsPBSF.Range("D30").Select() | Out-Null

# --- YS: row 23, selection moves to C23 ---------------------------------
$wsYS = $wb.Worksheets.Item("YS")
$wsYS.Range("A23").Value = 44096.490972222222
$wsYS.Range("B23").Value = 491.1
$wsYS.Range("C23").Value = 19.5
$wsYS.Range("C23").Select() | Out-Null

# --- SW: row 22, selection moves to C29 ---------------------------------
$wsSW = $wb.Worksheets.Item("SW")
$wsSW.Range("A22").Value = 44096.511111111111
$wsSW.Range("B22").Value = 1064
$wsSW.Range("C22").Value = 17.5
$wsSW.Range("C29").Select() | Out-Null

# --- YI: row 21, selection moves to C23 ---------------------------------
$wsYI = $wb.Worksheets.Item("YI")
$wsYI.Range("A21").Value = 44096.402777777781
$wsYI.Range("B21").Value = 442.5
$wsYI.Range("C21").Value = 18.5
$wsYI.Range("C23").Select() | Out-Null

# --- YN: row 21, selection moves to C21 ---------------------------------
$wsYN = $wb.Worksheets.Item("YN")
$wsYN.Range("A21").Value = 44096.418749999997
$wsYN.Range("B21").Value = 571.1
$wsYN.Range("C21").Value = 18.100000000000001
$wsYN.Range("C21").Select() | Out-Null

# --- 6MC: row 22, selection moves to B25 --------------------------------
$ws6MC = $wb.Worksheets.Item("6MC")
$ws6MC.Range("A22").Value = 44096.436111111114
$ws6MC.Range("B22").Value = 674
$ws6MC.Range("C22").Value = 15.5
$ws6MC.Range("B25").Select() | Out-Null

# --- DC: row 22, selection moves to H18 ---------------------------------
$wsDC = $wb.Worksheets.Item("DC")
$wsDC.Range("A22").Value = 44096.444444444445
$wsDC.Range("B22").Value = 679.2
$wsDC.Range("C22").Value = 14.5
$wsDC.Range("H18").Select() | Out-Null

# --- PBMS: row 23, selection moves to D34 -------------------------------
$wsPBMS = $wb.Worksheets.Item("PBMS")
$wsPBMS.Range("A23").Value = 44096.459722222222
$wsPBMS.Range("B23").Value = 980.5
$wsPBMS.Range("C23").Value = 16.899999999999999
$wsPBMS.Range("D34").Select() | Out-Null

# --- WIC: row 8, becomes the active tab, selection moves to H14 --------
$wsWIC = $wb.Worksheets.Item("WIC")
$wsWIC.Range("A8").Value = 44096.388888888891
$wsWIC.Range("B8").Value = 682.7
$wsWIC.Range("C8").Value = 17.7
$wsWIC.Activate()
$wsWIC.Range("H14").Select() | Out-Null
